$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "ujang"
$ws.Range("A2").Value = "imas"
$ws.Range("A3").Value = "ari"

$ws.Range("A4").Select()
